# Adapt column header formatting to respective input file names.
#
# The sheet compares an "old" (FV2210) AHB file against a "new" (FV2304)
# AHB file. The header row used generic "_old" / "_new" suffixes; rename
# them to "_FV2210" / "_FV2304" respectively (columns A:J and L:U - the
# "diff" column K stays as-is). Then freeze the header row and wrap the
# data range in a real Excel Table (ListObject) with AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells -------------------------------------------------
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A-J (1-10) -> "_old" becomes "_FV2210"
for ($i = 0; $i -lt $fv2210Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}

# Column K (11) is "diff" and is left untouched.

# Columns L-U (12-21) -> "_new" becomes "_FV2304"
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# --- 2. Freeze the header row ----------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into an Excel Table (adds AutoFilter too) -----------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = $null
